# Nexial showcase workbook update: add the "aws.ses" command category
# (sendMail / sendTextMail) to the hidden '#system' lookup sheet, and
# register the new named range. This mirrors the commit that introduced
# AWS SES mailer support to Nexial.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# ---------------------------------------------------------------------
# 1) Make room: insert a new column before column C. This shifts the
#    existing per-category command columns C:Z to D:AA, leaving column
#    A (the "target" index) and column B (aws.s3) untouched.
# ---------------------------------------------------------------------
$ws.Columns("C:C").Insert()

# ---------------------------------------------------------------------
# 2) The "target" index in column A lists every category name, sorted
#    alphabetically. Insert "aws.ses" between "aws.s3" and "base",
#    pushing the remaining entries down by one row.
# ---------------------------------------------------------------------
$ws.Range("A27").Value = "xml"
$ws.Range("A26").Value = "ws.async"
$ws.Range("A25").Value = "ws"
$ws.Range("A24").Value = "webcookie"
$ws.Range("A23").Value = "webalert"
$ws.Range("A22").Value = "web"
$ws.Range("A21").Value = "step"
$ws.Range("A20").Value = "ssh"
$ws.Range("A19").Value = "sound"
$ws.Range("A18").Value = "sms"
$ws.Range("A17").Value = "redis"
$ws.Range("A16").Value = "rdbms"
$ws.Range("A15").Value = "pdf"
$ws.Range("A14").Value = "number"
$ws.Range("A13").Value = "mail"
$ws.Range("A12").Value = "json"
$ws.Range("A11").Value = "jms"
$ws.Range("A10").Value = "io"
$ws.Range("A9").Value = "image"
$ws.Range("A8").Value = "external"
$ws.Range("A7").Value = "excel"
$ws.Range("A6").Value = "desktop"
$ws.Range("A5").Value = "csv"
$ws.Range("A4").Value = "base"
$ws.Range("A3").Value = "aws.ses"

# ---------------------------------------------------------------------
# 3) Populate the newly-inserted column C with the aws.ses category:
#    header + the two new commands.
# ---------------------------------------------------------------------
$ws.Range("C1").Value = "aws.ses"
$ws.Range("C2").Value = "sendMail(profile,to,subject,body)"
$ws.Range("C3").Value = "sendTextMail(profile,to,subject,body)"

# ---------------------------------------------------------------------
# 4) Re-point the existing named ranges at their new columns (C:Z -> D:AA)
#    and register the new "aws.ses" named range.
# ---------------------------------------------------------------------
$names = $wb.Names
$names.Item("base").RefersTo = "='#system'!`$D`$2:`$D`$36"
$names.Item("csv").RefersTo = "='#system'!`$E`$2:`$E`$5"
$names.Item("desktop").RefersTo = "='#system'!`$F`$2:`$F`$92"
$names.Item("excel").RefersTo = "='#system'!`$G`$2:`$G`$14"
$names.Item("external").RefersTo = "='#system'!`$H`$2:`$H`$3"
$names.Item("image").RefersTo = "='#system'!`$I`$2:`$I`$5"
$names.Item("io").RefersTo = "='#system'!`$J`$2:`$J`$24"
$names.Item("jms").RefersTo = "='#system'!`$K`$2:`$K`$4"
$names.Item("json").RefersTo = "='#system'!`$L`$2:`$L`$14"
$names.Item("mail").RefersTo = "='#system'!`$M`$2:`$M`$2"
$names.Item("number").RefersTo = "='#system'!`$N`$2:`$N`$15"
$names.Item("pdf").RefersTo = "='#system'!`$O`$2:`$O`$16"
$names.Item("rdbms").RefersTo = "='#system'!`$P`$2:`$P`$7"
$names.Item("redis").RefersTo = "='#system'!`$Q`$2:`$Q`$10"
$names.Item("sms").RefersTo = "='#system'!`$R`$2:`$R`$2"
$names.Item("sound").RefersTo = "='#system'!`$S`$2:`$S`$5"
$names.Item("ssh").RefersTo = "='#system'!`$T`$2:`$T`$9"
$names.Item("step").RefersTo = "='#system'!`$U`$2:`$U`$4"
$names.Item("target").RefersTo = "='#system'!`$A`$2:`$A`$27"
$names.Item("web").RefersTo = "='#system'!`$V`$2:`$V`$117"
$names.Item("webalert").RefersTo = "='#system'!`$W`$2:`$W`$8"
$names.Item("webcookie").RefersTo = "='#system'!`$X`$2:`$X`$8"
$names.Item("ws").RefersTo = "='#system'!`$Y`$2:`$Y`$17"
$names.Item("ws.async").RefersTo = "='#system'!`$Z`$2:`$Z`$8"
$names.Item("xml").RefersTo = "='#system'!`$AA`$2:`$AA`$11"

$names.Add("aws.ses", "='#system'!`$C`$2:`$C`$3")

Write-Host "aws.ses command category added."
